$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOC")
$ws.Activate()
$ws.Range("B2").Value = 15
$ws.Range("B5").Select()
